$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the start/end date cells (Fecha de inicio / Fecha de fin).
# Values are written as Excel date serials (days since 1899-12-30) so no
# implicit/default date style gets auto-created by a DateTime assignment.
$ws.Range("D2").Value = 38772
$ws.Range("E2").Value = 45716

$ws.Range("D3").Value = 38773
$ws.Range("E3").Value = 46021

$ws.Range("D4").Value = 38774
$ws.Range("E4").Value = 46368

$ws.Range("D5").Value = 38775
$ws.Range("E5").Value = 47832

# Apply the date number format (built-in format id 14, "mm-dd-yy") to D2 first,
# then propagate that exact style to the rest of the range via PasteSpecial so
# every touched cell shares a single style index.
$ws.Range("D2").NumberFormat = "mm-dd-yy"
$ws.Range("D2").Copy()
$ws.Range("D2:E5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the saved selection on the sheet
$ws.Range("F6").Select()
